# "Generate Report for handoff"
#
# The localization-status report is regenerated: the row describing the
# 09cf1117-...md file moves from "Handed back" into a fresh handoff round
# ("Ready for handoff", with a new Latest Handoff Datetime), while the
# 4ed2057b-...md row keeps its "Handed back: in sync with en-US" status.
# On the Overview sheet (and as a side-effect of the regeneration on the
# zh-cn / de-de detail sheets) the two file rows end up swapped: row 2 now
# describes 4ed2057b-...md and row 3 now describes 09cf1117-...md. The
# existing hyperlink objects keep pointing at their original target URLs
# (same r:id), only their displayed text moves together with the cell
# text.

$wb = $excel.ActiveWorkbook

function To-AbsAddr {
    param($addr)
    if ($addr -match '^([A-Za-z]+)(\d+)$') {
        return '$' + $matches[1] + '$' + $matches[2]
    }
    return $addr
}

# Sets a cell's value and, if a hyperlink is anchored to that cell, keeps
# the hyperlink's displayed text in sync (without touching its target URL
# / r:id).
function Set-CellText {
    param($ws, $addr, $value)
    $ws.Range($addr).Value = $value
    $absAddr = To-AbsAddr $addr
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $absAddr) {
            $h.TextToDisplay = $value
        }
    }
}

$uuid09 = "09cf1117-62f3-4130-aa59-b00a50001643.md"
$uuid4e = "4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.md"

$handedBack = "Handed back: in sync with en-US"
$readyForHandoff = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: rows 2/3 swap which file they describe.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

Set-CellText $ws "A2" $uuid4e
Set-CellText $ws "B2" $handedBack
Set-CellText $ws "C2" $handedBack

Set-CellText $ws "A3" $uuid09
Set-CellText $ws "B3" $readyForHandoff
Set-CellText $ws "C3" $readyForHandoff

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$xlf4eZh = "4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.2ae90447bf1701606f56496466a12aeea19087ff.zh-cn.xlf"
$xlf09Zh = "09cf1117-62f3-4130-aa59-b00a50001643.6948b8e253bdeb612cb6e16789f274eeb6ea7b25.zh-cn.xlf"

Set-CellText $ws "A2" $uuid4e
Set-CellText $ws "B2" $handedBack
Set-CellText $ws "C2" $xlf4eZh
Set-CellText $ws "D2" "2016-01-26 05:36:14"
Set-CellText $ws "E2" $uuid4e
Set-CellText $ws "F2" $xlf4eZh
Set-CellText $ws "G2" "2016-01-26 05:37:07"
Set-CellText $ws "H2" "Include"

Set-CellText $ws "A3" $uuid09
Set-CellText $ws "B3" $readyForHandoff
Set-CellText $ws "C3" $xlf09Zh
Set-CellText $ws "D3" "2016-01-26 05:38:08"
Set-CellText $ws "E3" $uuid09
Set-CellText $ws "F3" $xlf09Zh
Set-CellText $ws "G3" "2016-01-26 05:37:07"
Set-CellText $ws "H3" "Include"

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$xlf4eDe = "4ed2057b-0b97-4a12-8431-67f2d9a1b8d4.2ae90447bf1701606f56496466a12aeea19087ff.de-de.xlf"
$xlf09De = "09cf1117-62f3-4130-aa59-b00a50001643.6948b8e253bdeb612cb6e16789f274eeb6ea7b25.de-de.xlf"

Set-CellText $ws "A2" $uuid4e
Set-CellText $ws "B2" $handedBack
Set-CellText $ws "C2" $xlf4eDe
Set-CellText $ws "D2" "2016-01-26 05:36:26"
Set-CellText $ws "E2" $uuid4e
Set-CellText $ws "F2" $xlf4eDe
Set-CellText $ws "G2" "2016-01-26 05:37:24"
Set-CellText $ws "H2" "Include"

Set-CellText $ws "A3" $uuid09
Set-CellText $ws "B3" $readyForHandoff
Set-CellText $ws "C3" $xlf09De
Set-CellText $ws "D3" "2016-01-26 05:38:23"
Set-CellText $ws "E3" $uuid09
Set-CellText $ws "F3" $xlf09De
Set-CellText $ws "G3" "2016-01-26 05:37:24"
Set-CellText $ws "H3" "Include"
